$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.827.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.363.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.48%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.361.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.04%  "

$ws.Range("E11").Value = "  +2.43%  "

$ws.Range("E12").Value = "  +3.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.935.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.363.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.999.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.22%  "

$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  +2.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.496.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("E27").Value = "  +11.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +21.90%  "

$ws.Range("E29").Value = "  +11.35%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.61%  "

$ws.Range("E32").Value = "  +2.11%  "

$ws.Range("E33").Value = "  +4.69%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.393.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  +3.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.56%  "

$ws.Range("E38").Value = "  +4.82%  "

$ws.Range("E39").Value = "  +5.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.76%  "

$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.40%  "

$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.91%  "

$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("E47").Value = "  +4.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.26%  "
